$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1886
$ws1.Range("F5").Value = 20
$ws1.Range("F6").Value = 854
$ws1.Range("F9").Value = 42
$ws1.Range("F11").Value = 247
$ws1.Range("F12").Value = 17
$ws1.Range("F14").Value = 141
$ws1.Range("F16").Value = 4444
$ws1.Range("F19").Value = 483
$ws1.Range("F20").Value = 433
$ws1.Range("F24").Value = 2085
$ws1.Range("F26").Value = 52
$ws1.Range("F27").Value = 33
$ws1.Range("F28").Value = 52
$ws1.Range("F29").Value = 2134
$ws1.Range("F30").Value = 79
$ws1.Range("F35").Value = 37
$ws1.Range("F37").Value = 31

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1886
$ws4.Range("F5").Value = 20
$ws4.Range("F6").Value = 854
$ws4.Range("F9").Value = 42
$ws4.Range("F11").Value = 247
$ws4.Range("F12").Value = 17
$ws4.Range("F14").Value = 141
$ws4.Range("F17").Value = 4444
$ws4.Range("F20").Value = 483
$ws4.Range("F21").Value = 433
$ws4.Range("F25").Value = 2085
$ws4.Range("F27").Value = 52
$ws4.Range("F28").Value = 33
$ws4.Range("F29").Value = 52
$ws4.Range("F30").Value = 2134
$ws4.Range("F31").Value = 79
$ws4.Range("F36").Value = 37
$ws4.Range("F38").Value = 31
